# Apply the "thay đôi chiến lược chạy multi process. Sửa lại template báo
# cáo tổng hợp cơ sở" edit to the workbook.
#
# Sheet 2 "CHI TIẾT VỀ THU NỢ": insert a new "Ngày thực hiện" column before
#   the existing "Lượng thu" column (F -> G).
# Sheet 7 "QUỸ LƯƠNG": update several employees' computed salary totals
#   (and the grand total) to reflect the new multi-process run.
# Sheet 8 "LỢI NHUẬN": replace the old single-base summary columns with the
#   new, wider "tổng hợp cơ sở" template (Cơ sở / Tổng đơn giá / Đã thanh
#   toán / ... / Tỉ lệ lợi nhuận) and fill in the recomputed row for
#   LONG XUYÊN.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "CHI TIẾT VỀ THU NỢ": insert "Ngày thực hiện" column before "Lượng thu"
# ---------------------------------------------------------------------
$wsThuNo = $wb.Worksheets.Item(2)

# Shifts the old column F ("Lượng thu") one place to the right, to G,
# and opens up a blank column F.
$wsThuNo.Columns.Item(6).Insert()

$wsThuNo.Cells.Item(1, 6).Value = "Ngày thực hiện"
$wsThuNo.Cells.Item(2, 6).Value = ""

# ---------------------------------------------------------------------
# Sheet "QUỸ LƯƠNG": recomputed salary totals after the strategy change
# ---------------------------------------------------------------------
$wsQuyLuong = $wb.Worksheets.Item(7)

$wsQuyLuong.Cells.Item(10, 3).Value = 53571.42857142857
$wsQuyLuong.Cells.Item(12, 3).Value = 1957142.857142857
$wsQuyLuong.Cells.Item(13, 3).Value = 535714.2857142857
$wsQuyLuong.Cells.Item(14, 3).Value = 528571.4285714286
$wsQuyLuong.Cells.Item(15, 3).Value = 640714.2857142857
$wsQuyLuong.Cells.Item(16, 3).Value = 355357.1428571428
$wsQuyLuong.Cells.Item(22, 3).Value = 6623452.380952381

# ---------------------------------------------------------------------
# Sheet "LỢI NHUẬN": new wider "báo cáo tổng hợp cơ sở" template
# ---------------------------------------------------------------------
$wsLoiNhuan = $wb.Worksheets.Item(8)

$headers = @("Cơ sở", "Tổng đơn giá", "Đã thanh toán", "Tỉ lệ thanh toán", "Tỉ lệ nợ", "Thu nợ", "Tổng doanh thu", "Chi tiêu", "Quỹ lương", "Tổng chi phí", "Lợi nhuận", "Tỉ lệ lợi nhuận")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsLoiNhuan.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$wsLoiNhuan.Cells.Item(2, 1).Value = "LONG XUYÊN"
$wsLoiNhuan.Cells.Item(2, 2).Value = 6000000
$wsLoiNhuan.Cells.Item(2, 3).Value = 5500000
$wsLoiNhuan.Cells.Item(2, 4).Value = 0.9166666666666666
$wsLoiNhuan.Cells.Item(2, 5).Value = 0.08333333333333337
$wsLoiNhuan.Cells.Item(2, 6).Value = 0
$wsLoiNhuan.Cells.Item(2, 7).Value = 5500000
$wsLoiNhuan.Cells.Item(2, 8).Value = 928500
$wsLoiNhuan.Cells.Item(2, 9).Value = 6623452.380952381
$wsLoiNhuan.Cells.Item(2, 10).Value = 7551952.380952381
$wsLoiNhuan.Cells.Item(2, 11).Value = -2051952.380952381
$wsLoiNhuan.Cells.Item(2, 12).Value = -0.373082251082251

Write-Output "edit applied"
